$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-05 01:22:36"

for ($row = 2; $row -le 21; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
